# feat: add 2022-Q3 data
#
# 1. "总计" (summary) sheet: insert a new row right under the header for the
#    2022-Q3 quarter (holdings count = 2, holding value = 0), pushing every
#    existing quarter row down by one and re-numbering the running index in
#    column A.
# 2. Insert a brand-new worksheet named "2022-Q3" right after "总计" (i.e.
#    before the existing "2022-Q2" sheet) holding the two funds reported for
#    that quarter.

$wb = $excel.ActiveWorkbook

function Set-BoldBoxStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous
}

function Set-TextValue($cell, [string]$text) {
    # Force text storage so numeric-looking strings ("0.26", "005502", ...)
    # are not silently reinterpreted as numbers (and don't lose leading
    # zeros), matching the source data which stores these as plain text.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert the 2022-Q3 summary row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
Set-BoldBoxStyle($total.Range("A2"))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0

# Rows 3..9 now hold what used to be rows 2..8; refresh the running index in
# column A (0-based) to stay contiguous with the inserted row.
for ($r = 3; $r -le 9; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# Step 2: new "2022-Q3" sheet, positioned right after "总计"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}
Set-BoldBoxStyle($newSheet.Range("B1:H1"))

# Row 2: 005502 华泰紫金智能量化股票A
Set-BoldBoxStyle($newSheet.Range("A2"))
$newSheet.Cells.Item(2, 1).Value = 0
Set-TextValue $newSheet.Cells.Item(2, 2) "005502"
$newSheet.Cells.Item(2, 3).Value = "华泰紫金智能量化股票A"
Set-TextValue $newSheet.Cells.Item(2, 4) "0.26"
Set-TextValue $newSheet.Cells.Item(2, 5) "94.33"
Set-TextValue $newSheet.Cells.Item(2, 6) "1.09"
Set-TextValue $newSheet.Cells.Item(2, 7) "0.0028"
$newSheet.Cells.Item(2, 8).Value = 10

# Row 3: 014629 华泰紫金智能量化股票C
Set-BoldBoxStyle($newSheet.Range("A3"))
$newSheet.Cells.Item(3, 1).Value = 1
Set-TextValue $newSheet.Cells.Item(3, 2) "014629"
$newSheet.Cells.Item(3, 3).Value = "华泰紫金智能量化股票C"
Set-TextValue $newSheet.Cells.Item(3, 4) "0.00"
Set-TextValue $newSheet.Cells.Item(3, 5) "94.33"
Set-TextValue $newSheet.Cells.Item(3, 6) "1.09"
$newSheet.Cells.Item(3, 7).Value = 0
$newSheet.Cells.Item(3, 8).Value = 10
